# Commit: "alpha depends on 2d/3d filter"
#
# Sheet CHA_RETAU550 (xl/worksheets/sheet4.xml): the old row 34
# (NOSLIP+DSMAG+SMALL, grid 96x64x64) is replaced by three rows:
#   - row 35: same case family but on the finer 96x64x128 grid (new alpha
#             numbers recomputed for that grid -- this used to be an
#             untested grid for the dynamic-Smagorinsky case)
#   - row 36: the original 96x64x64 numbers, re-added unchanged
#   - row 37: a brand new variant ("...+UPDATE1") on the 96x64x64 grid,
#             with an alpha value that now differs slightly because the
#             filter is computed as 2d/3d depending on the test case.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CHA_RETAU550")

# --- Remove the old row 34; rows 35+ will be (re)written from scratch,
#     leaving a genuine gap at row 34 (matches the target file, which has
#     no <row r="34"> element at all). ---
$ws.Rows.Item(34).Delete()

# Carry over the per-column number formatting (styles) that the row used
# to have, since the new rows are being created from nothing. Row 33 uses
# exactly the same per-column styles as the old row 34 did.
$ws.Range("A33:X33").Copy()
$ws.Range("A35:X37").PasteSpecial(-4122)

# ============================== Row 35 ==============================
# NOSLIP+DSMAG+SMALL on the  96x64x128 grid.
# NB: numeric literals are written in plain decimal (not scientific
# notation) -- the COM script parser here does not accept "1.23E-4" style
# literals.
$ws.Range("A35").Value = "NOSLIP+DSMAG+SMALL"
$ws.Range("B35").Value = 20540
$ws.Range("C35").Value = "6.4×2.4×2.0"
$ws.Range("D35").Value = " 96×64×128"
$ws.Range("E35").Value = 0.0666666666666666
$ws.Range("F35").Value = 0.0375
$ws.Range("G35").Value = 0.03957235
$ws.Range("H35").Value = 0.001094476
$ws.Range("I35").Formula = "=E35/F35"
$ws.Range("J35").Value = "(gtype,gr)=(1,5)"
$ws.Range("K35").Value = "CFR"
$ws.Range("L35").Value = "NA"
$ws.Range("M35").Value = "NA"
$ws.Range("N35").Value = 551.75
$ws.Range("O35").Value = 546.751854551584
$ws.Range("P35").Formula = "=8*(N35/B35)^2"
$ws.Range("Q35").Formula = "=8*(O35/B35)^2"
$ws.Range("R35").Formula = "=(Q35-P35)/P35"
$ws.Range("S35").Formula = "=200*2*O35/B35"
$ws.Range("T35").Formula = "=B35/4*P35"
$ws.Range("U35").Formula = "=E35*N35"
$ws.Range("V35").Formula = "=F35*N35"
$ws.Range("W35").Formula = "=G35*N35"
$ws.Range("X35").Formula = "=H35*N35"

# ============================== Row 36 ==============================
# NOSLIP+DSMAG+SMALL on the original 96x64x64 grid (same values the old
# row 34 had).
$ws.Range("A36").Value = "NOSLIP+DSMAG+SMALL"
$ws.Range("B36").Value = 20540
$ws.Range("C36").Value = "6.4×2.4×2.0"
$ws.Range("D36").Value = " 96×64×64"
$ws.Range("E36").Value = 0.0666666666666666
$ws.Range("F36").Value = 0.0375
$ws.Range("G36").Value = 0.07902424
$ws.Range("H36").Value = 0.002276549
$ws.Range("I36").Formula = "=E36/F36"
$ws.Range("J36").Value = "(gtype,gr)=(1,5)"
$ws.Range("K36").Value = "CFR"
$ws.Range("L36").Value = "NA"
$ws.Range("M36").Value = "NA"
$ws.Range("N36").Value = 551.75
$ws.Range("O36").Value = 530.154008126688
$ws.Range("P36").Formula = "=8*(N36/B36)^2"
$ws.Range("Q36").Formula = "=8*(O36/B36)^2"
$ws.Range("R36").Formula = "=(Q36-P36)/P36"
$ws.Range("S36").Value = 11
$ws.Range("T36").Formula = "=B36/4*P36"
$ws.Range("U36").Formula = "=E36*N36"
$ws.Range("V36").Formula = "=F36*N36"
$ws.Range("W36").Formula = "=G36*N36"
$ws.Range("X36").Formula = "=H36*N36"

# ============================== Row 37 ==============================
# New case: NOSLIP+DSMAG+SMALL+UPDATE1 on the 96x64x64 grid -- alpha now
# depends on whether the filter is 2d or 3d, giving a (slightly) different
# O37 value than row 36's O36.
$ws.Range("A37").Value = "NOSLIP+DSMAG+SMALL+UPDATE1"
$ws.Range("B37").Value = 20540
$ws.Range("C37").Value = "6.4×2.4×2.0"
$ws.Range("D37").Value = " 96×64×64"
$ws.Range("E37").Value = 0.0666666666666666
$ws.Range("F37").Value = 0.0375
$ws.Range("G37").Value = 0.07902424
$ws.Range("H37").Value = 0.002276549
$ws.Range("I37").Formula = "=E37/F37"
$ws.Range("J37").Value = "(gtype,gr)=(1,5)"
$ws.Range("K37").Value = "CFR"
$ws.Range("L37").Value = "NA"
$ws.Range("M37").Value = "NA"
$ws.Range("N37").Value = 551.75
$ws.Range("O37").Value = 530.109409358885
$ws.Range("P37").Formula = "=8*(N37/B37)^2"
$ws.Range("Q37").Formula = "=8*(O37/B37)^2"
$ws.Range("R37").Formula = "=(Q37-P37)/P37"
$ws.Range("S37").Value = 11
$ws.Range("T37").Formula = "=B37/4*P37"
$ws.Range("U37").Formula = "=E37*N37"
$ws.Range("V37").Formula = "=F37*N37"
$ws.Range("W37").Formula = "=G37*N37"
$ws.Range("X37").Formula = "=H37*N37"

# Leave the selection where the author left it.
$ws.Range("T42").Select()
